$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.774.79"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "2.078.84"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.67"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.84"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.774"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").Value = "2.079.21"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").Value = "37.701.80"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.135"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.59%  "

$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0635"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("E34").Value = "  -4.68%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -2.38%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0977"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.52%  "

$ws.Range("E41").Value = "  -2.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.07%  "

$ws.Range("D44").Value = "1.439.58"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("E49").Value = "  -2.07%  "

$ws.Range("D50").Value = "2.268.57"
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.12%  "
